# Add documentation row for the new "array" command to the commands sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 7 ("settings(...)" row), which shifts
# that row (and everything below it, including the merged description block)
# down by one - matching the target layout where the new "array(...)" entry
# becomes row 7 and the old rows 7/8 become rows 8/9.
$ws.Rows.Item(7).Insert()

# Copy the formatting used by the other command rows (row 4, "stroke(...)")
# onto the freshly inserted row so borders/alignment match the rest of the
# table.
$ws.Range("A4:C4").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new command's documentation.
$ws.Range("A7").Value2 = "array([n,m,o],color)"
$ws.Range("B7").Value2 = "array([2,6,9],0xFFFFFF)"
$ws.Range("C7").Value2 = "Sets multiple pixels in array a certain color "

# Match the saved selection state recorded in the workbook.
$ws.Range("C7").Select()

# The workbook was also re-saved with explicit page setup (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
